$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet/tab to reflect the new "through" date
$ws.Name = "Through 2021-10-10"

# Row 12 - October (through 10-10) data
$ws.Range("A12").Value = "October (through 10-10)"
$ws.Range("C12").Value = 7
$ws.Range("D12").Value = 0.125
$ws.Range("F12").Value = 17
$ws.Range("L12").Value = 23
$ws.Range("M12").Value = 0.08
$ws.Range("O12").Value = 11
$ws.Range("R12").Value = 43
$ws.Range("U12").Value = 66

# Row 13 - Total data
$ws.Range("C13").Value = 203
$ws.Range("D13").Value = 0.1325
$ws.Range("F13").Value = 400
$ws.Range("G13").Value = 0.1031
$ws.Range("L13").Value = 510
$ws.Range("M13").Value = 0.1099
$ws.Range("O13").Value = 390
$ws.Range("P13").Value = 0.0993
$ws.Range("R13").Value = 891
$ws.Range("S13").Value = 0.0561
$ws.Range("U13").Value = 1237
$ws.Range("V13").Value = 0.06
